$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$rows = @(
    @('17:42:06', 'In Bed | HR=49 | BR=1'),
    @('17:42:06', 'In Bed | HR=50 | BR=2'),
    @('17:42:06', 'In Bed | HR=53 | BR=5'),
    @('17:42:07', 'In Bed | HR=50 | BR=2'),
    @('17:42:07', 'In Bed | HR=55 | BR=7'),
    @('17:42:07', 'In Bed | HR=57 | BR=9'),
    @('17:42:08', 'In Bed | HR=77 | BR=29'),
    @('17:42:10', 'In Bed | HR=71 | BR=23'),
    @('17:42:10', 'In Bed | HR=73 | BR=25'),
    @('17:42:11', 'In Bed | HR=76 | BR=28'),
    @('17:42:12', 'In Bed | HR=56 | BR=8'),
    @('17:42:13', 'In Bed | HR=54 | BR=6'),
    @('17:42:14', 'In Bed | HR=75 | BR=27'),
    @('17:42:15', 'In Bed | HR=52 | BR=4'),
    @('17:42:16', 'In Bed | HR=54 | BR=6'),
    @('17:42:17', 'In Bed | HR=50 | BR=2'),
    @('17:42:19', 'In Bed | HR=98 | BR=50'),
    @('17:42:20', 'In Bed | HR=68 | BR=20'),
    @('17:42:21', 'In Bed | HR=89 | BR=41'),
    @('17:42:22', 'In Bed | HR=105 | BR=57'),
    @('17:42:23', 'In Bed | HR=56 | BR=8'),
    @('17:42:24', 'In Bed | HR=54 | BR=6'),
    @('17:42:25', 'In Bed | HR=105 | BR=57'),
    @('17:42:27', 'In Bed | HR=103 | BR=55'),
    @('17:42:27', 'In Bed | HR=94 | BR=46'),
    @('17:42:28', 'In Bed | HR=66 | BR=18'),
    @('17:42:29', 'In Bed | HR=101 | BR=53'),
    @('17:42:30', 'In Bed | HR=75 | BR=27'),
    @('17:42:31', 'In Bed | HR=52 | BR=4'),
    @('17:42:32', 'In Bed | HR=50 | BR=2'),
    @('17:42:33', 'In Bed | HR=54 | BR=6'),
    @('17:42:34', 'In Bed | HR=50 | BR=2'),
    @('17:42:43', 'In Bed | HR=51 | BR=3'),
    @('17:42:44', 'In Bed | HR=50 | BR=2')
)

$startRow = 149
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $timestamp = $rows[$i][0]
    $value = $rows[$i][1]
    $ws.Cells.Item($r, 1).Value = "'2026-01-28"
    $ws.Cells.Item($r, 2).Value = $timestamp
    $ws.Cells.Item($r, 3).Value = "17:00"
    $ws.Cells.Item($r, 4).Value = "Bedroom"
    $ws.Cells.Item($r, 5).Value = $value
    $ws.Cells.Item($r, 6).Value = "Occupied"
}
